# Daily attendance processing - 2025-11-26 07:26:56
# Normalizes the "Recorded By" (column G) comma-separated list of
# recorders on each attendance row so that the literal entry "System"
# always appears first, followed by the remaining recorders sorted
# alphabetically (case-insensitive).

function Compare-RecordedByPart($x, $y) {
    $xIsSystem = $x.Equals("System")
    $yIsSystem = $y.Equals("System")
    if ($xIsSystem -and -not $yIsSystem) { return -1 }
    if ($yIsSystem -and -not $xIsSystem) { return 1 }
    if ($xIsSystem -and $yIsSystem) { return 0 }

    $xl = $x.ToLower()
    $yl = $y.ToLower()
    if ($xl -lt $yl) { return -1 }
    if ($xl -gt $yl) { return 1 }
    return 0
}

function Sort-RecordedBy($value) {
    $parts = $value -split ", "
    $n = $parts.Length

    for ($i = 1; $i -lt $n; $i++) {
        $key = $parts[$i]
        $j = $i - 1
        while ($j -ge 0 -and (Compare-RecordedByPart $parts[$j] $key) -gt 0) {
            $parts[$j + 1] = $parts[$j]
            $j = $j - 1
        }
        $parts[$j + 1] = $key
    }

    return ($parts -join ", ")
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Value2
    if ($current -ne $null -and $current -ne "") {
        $updated = Sort-RecordedBy $current
        if ($updated -ne $current) {
            $cell.Value2 = $updated
        }
    }
}
